{"js": "/*\n * Edit: In the TCP Session Hijacking task paragraph, merge the trailing\n * \"to perform this task.\" / \"Use the ...hexify.py... ascii text.\"\n * runs, append a separating space, and add new guidance sentences about\n * providing psh/ack flags and the hijacking goal (deleting\n * ~/documents/delete-this.txt), matching the author's commit\n * \"give hints on use of nping parameters for spoofing\".\n */\nconst body = context.document.body;\n\n// Locate the unique run of text that ends the target paragraph.\nconst searchResults = body.search('Use the \u2013data option to send your payload.  Your attacker home directory includes a \u201chexify.py\u201d script that creates hex versions of ascii text.', { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find target sentence to edit.\");\n}\n\nconst targetRange = searchResults.items[0];\n\n// Append the new guidance text right after the existing sentence; it\n// inherits the surrounding run's formatting (Times, 11pt, black).\ntargetRange.insertText(' You will also want to provide the psh and ack flags, and ack the previous packet in your spoofed packet. Your goal is to use a spoofed packet to hijack a telnet session and delete the file on the server at ~/documents/delete-this.txt.  Note that if you use your telnet session to delete that file, e.g., to observe the protocol in wireshark, then you must recreate that file so it can be deleted in a hijacked session.', \"End\");\n\nawait context.sync();\n", "ps1": "# Edit: In the TCP Session Hijacking task paragraph, merge the trailing\n# \"to perform this task.\" / \"Use the ...hexify.py... ascii text.\" runs,\n# append a separating space, and add new guidance sentences about\n# providing psh/ack flags and the hijacking goal (deleting\n# ~/documents/delete-this.txt), matching the author's commit\n# \"give hints on use of nping parameters for spoofing\".\n\n$d = $word.ActiveDocument\n\n# Locate the unique run of text that ends the target paragraph.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Use the \u2013data option to send your payload.  Your attacker home directory includes a \u201chexify.py\u201d script that creates hex versions of ascii text.\"\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find target sentence to edit.\"\n}\n\n# Collapse the found range to its end point, then insert the new\n# guidance text there so it inherits the preceding run's formatting\n# (Times, 11pt, black).\n$rng.Collapse(0)\n$rng.InsertAfter(\" You will also want to provide the psh and ack flags, and ack the previous packet in your spoofed packet. Your goal is to use a spoofed packet to hijack a telnet session and delete the file on the server at ~/documents/delete-this.txt.  Note that if you use your telnet session to delete that file, e.g., to observe the protocol in wireshark, then you must recreate that file so it can be deleted in a hijacked session.\")\n$rng.Font.Name = \"Times\"\n$rng.Font.Size = 11\n$rng.Font.Color = 0\n"}
